$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2025-07-23 12:54:16"
$ws.Range("B7").Value = "delete-team"
$ws.Range("C7").Value = "new-organization97"
$ws.Range("D7").Value = "secondteam"
$ws.Range("E7").Value = "demo"
$ws.Range("I7").Value = "'False"
